$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC_ViewCart")

$ws.Cells.Item(2, 12).Value = 'VerifyTitle: City Market Norwalk - Online Grocery Supermarket with Home Delivery'
$ws.Cells.Item(3, 12).Value = 'VerifyElement: null'
$ws.Cells.Item(4, 12).Value = 'Click: null'
$ws.Cells.Item(5, 12).Value = 'VerifyElement: null'
$ws.Cells.Item(6, 12).Value = 'SetText: Randomemailid'
$ws.Cells.Item(7, 12).Value = 'SetText: 123456'
$ws.Cells.Item(8, 12).Value = 'Click: null'
$ws.Cells.Item(9, 12).Value = 'VerifyText: Akash sangal'
$ws.Cells.Item(10, 12).Value = 'VerifyTitle: City Market Norwalk - Online Grocery Supermarket with Home Delivery'
$ws.Cells.Item(11, 12).Value = 'VerifyElement: null'
$ws.Cells.Item(12, 12).Value = 'MoveToProductList: Quick & Easy Food Solutions'
$ws.Cells.Item(13, 12).Value = 'Click: null'
$ws.Cells.Item(14, 12).Value = 'VerifyElement: null'
$ws.Cells.Item(15, 12).Value = 'VerifyText: You have no items in your shopping cart.'
$ws.Cells.Item(16, 12).Value = 'Click: null'
$ws.Cells.Item(17, 12).Value = 'VerifyNoElement: null'
$ws.Cells.Item(18, 12).Value = 'MoveToProductList: Quick & Easy Food Solutions'
$ws.Cells.Item(19, 12).Value = 'MoveAndAddProduct: null'
$ws.Cells.Item(20, 12).Value = 'MoveAndAddProduct: exist'
$ws.Cells.Item(21, 12).Value = 'MoveAndAddProduct: Banquet Brown ''N Serve Turkey Sausage Links'
$ws.Cells.Item(22, 12).Value = 'MoveAndAddProduct: exist'
$ws.Cells.Item(23, 12).Value = 'MoveAndAddProduct: exist'
$ws.Cells.Item(24, 12).Value = 'MoveAndAddProduct: exist'
$ws.Cells.Item(25, 12).Value = 'MoveAndAddProduct: null'
$ws.Cells.Item(26, 12).Value = 'MoveAndAddProduct: exist'
$ws.Cells.Item(27, 12).Value = 'MoveAndAddProduct: exist'
$ws.Cells.Item(28, 12).Value = 'MoveAndAddProduct: null'
$ws.Cells.Item(29, 12).Value = 'MoveAndAddProduct: exist'
$ws.Cells.Item(30, 12).Value = 'Click: null'
$ws.Cells.Item(31, 12).Value = 'VerifyElement: null'
$ws.Cells.Item(32, 12).Value = 'Click: null'
$ws.Cells.Item(33, 12).Value = 'text not verified: Customers can opt for replacements in case an item is out of stock. We will choose a replacement item that is comparable to the original based on brand, flavor, size and price. You will be notified of replacements and will have the opportunity to approve the replacement or request a refund for the item.'
$ws.Cells.Item(33, 13).Value = 'Fail'
$ws.Cells.Item(33, 14).Value = 'Expected condition failed: waiting for presence of element located by: By.xpath: //div[normalize-space(@data-ui-id) = ''checkout-cart-validationmessages-message-error''] (tried for 40 second(s) with 500 MILLISECONDS interval)
Build info: version: ''unknown'', revision: ''1969d75'', time: ''2016-10-18 09:43:45 -0700''
System info: host: ''DESKTOP-OEL817D'', ip: ''192.168.134.2'', os.name: ''Windows 10'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''1.8.0_161''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities [{applicationCacheEnabled=false, rotatable=false, mobileEmulationEnabled=false, networkConnectionEnabled=false, chrome={chromedriverVersion=2.36.540470 (e522d04694c7ebea4ba8821272dbef4f9b818c91), userDataDir=C:\Users\Akash\AppData\Local\Temp\scoped_dir13564_28830}, takesHeapSnapshot=true, pageLoadStrategy=normal, databaseEnabled=false, handlesAlerts=true, hasTouchScreen=false, version=68.0.3440.106, platform=XP, browserConnectionEnabled=false, nativeEvents=true, acceptSslCerts=false, acceptInsecureCerts=false, locationContextEnabled=true, webStorageEnabled=true, browserName=chrome, takesScreenshot=true, javascriptEnabled=true, cssSelectorsEnabled=true, setWindowRect=true, unexpectedAlertBehaviour=}]
Session ID: 3fddb53966380820a83efadd593af767'
$ws.Cells.Item(34, 12).Value = 'AddInstruction: Hello'
$ws.Cells.Item(34, 13).Value = 'Pass'
$ws.Cells.Item(34, 14).Value = '-'
$ws.Cells.Item(35, 12).Value = 'AddInstruction: Hello'
$ws.Cells.Item(36, 12).Value = 'EditInstruction: Bye'
$ws.Cells.Item(37, 12).Value = 'RemoveInstruction: null'
$ws.Cells.Item(38, 12).Value = 'SetText: asdasd@'
$ws.Cells.Item(39, 12).Value = 'Wait: 6000'
$ws.Cells.Item(40, 12).Value = 'Click: null'
$ws.Cells.Item(41, 12).Value = 'VerifyText: Please enter a valid email address (Ex: johndoe@domain.com).'
$ws.Cells.Item(42, 12).Value = 'SetText: asdasd@assdcsadsaasd.com'
$ws.Cells.Item(43, 12).Value = 'Click: null'
$ws.Cells.Item(44, 12).Value = 'VerifyText: Thank you for your subscription.'
$ws.Cells.Item(45, 12).Value = 'SelectSubstitute: Do Not Allow Substitute'
$ws.Cells.Item(46, 12).Value = 'SelectSubstitute: Allow Substitute'
$ws.Cells.Item(47, 12).Value = 'SelectSubstitute: Allow Substitute'
$ws.Cells.Item(48, 12).Value = 'SelectSubstitute: Do Not Allow Substitute'
$ws.Cells.Item(49, 12).Value = 'SelectSubstitute: Allow Substitute'
$ws.Cells.Item(50, 12).Value = 'SelectSubstitute: Allow Substitute'
$ws.Cells.Item(51, 12).Value = 'VerifySummaryViewCart: null'
$ws.Cells.Item(52, 12).Value = 'AddProductfromViewCart: Jose Ole Steak & Cheese Chimichanga'
$ws.Cells.Item(53, 12).Value = 'VerifySummaryViewCart: null'
$ws.Cells.Item(53, 13).Value = 'Pass'
$ws.Cells.Item(53, 14).Value = '-'
$ws.Cells.Item(54, 12).Value = 'RemoveProductfromViewCart: Jose Ole Steak & Cheese Chimichanga'
$ws.Cells.Item(55, 12).Value = 'VerifySummaryViewCart: null'
$ws.Cells.Item(56, 12).Value = 'DeleteProductfromViewCart: Jose Ole Steak & Cheese Chimichanga'
$ws.Cells.Item(57, 12).Value = 'DeleteProductfromViewCart: all'
$ws.Cells.Item(58, 12).Value = 'VerifySummaryViewCart: null'
$ws.Cells.Item(59, 12).Value = 'VerifyFooterLinks: null'
$ws.Cells.Item(60, 12).Value = 'Click: null'
$ws.Cells.Item(61, 12).Value = 'Click: null'
$ws.Cells.Item(62, 12).Value = 'Wait: 6000'
$ws.Cells.Item(63, 12).Value = 'VerifyElement: null'
